$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(19).Delete()
$null = $ws.Range("C18").Select()
